$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3 values (was "Test" row, becomes "Fast" row)
$ws.Range("A3").Value = "Fast"
$ws.Range("B3").Value = 3
$ws.Range("C3").Value = 8
$ws.Range("D3").Value = 10
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 8
$ws.Range("G3").Value = 1

# Row 4 values (was "Fast" row, becomes "Test" row with new values)
$ws.Range("A4").Value = "Test"
$ws.Range("B4").Value = 5
$ws.Range("C4").Value = 0.4
$ws.Range("D4").Value = 44
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 20
$ws.Range("G4").Value = 1
